$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F7").Value = "93_referral_statement"
$ws.Range("F18").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F22").Value = "ppe || application instructions"
$ws.Range("F24").Value = "ppe"
$ws.Range("F31").Value = "env warning - water"
$ws.Range("F32").Value = "application instructions || env warning - species"
$ws.Range("F55").Value = "application instructions"
$ws.Range("F56").Value = "application instructions"
$ws.Range("F58").Value = "application instructions"
$ws.Range("F59").Value = "application instructions"
$ws.Range("F63").Value = "use restrictions"
$ws.Range("F64").Value = "use restrictions"
$ws.Range("F66").Value = "135_product_information"
$ws.Range("F68").Value = "application instructions"
$ws.Range("F69").Value = "mixing"
$ws.Range("F70").Value = "use restrictions"
$ws.Range("F72").Value = "application instructions"
$ws.Range("F86").Value = "use restrictions"
$ws.Range("F87").Value = "off target movement"
$ws.Range("F88").Value = "off target movement"
$ws.Range("F89").Value = "off target movement"
$ws.Range("F91").Value = "off target movement"
$ws.Range("F92").Value = "off target movement"
$ws.Range("F93").Value = "off target movement"
$ws.Range("F94").Value = "off target movement"
$ws.Range("F95").Value = "off target movement"
$ws.Range("F98").Value = "application instructions"
$ws.Range("F99").Value = "application instructions"
$ws.Range("F102").Value = "application instructions"
$ws.Range("F105").Value = "application instructions"
$ws.Range("F108").Value = "application instructions"
$ws.Range("F109").Value = "application instructions"
$ws.Range("F111").Value = "safety procedures || application instructions"
$ws.Range("F113").Value = "mixing"
$ws.Range("F114").Value = "mixing"
$ws.Range("F116").Value = "use restrictions"
$ws.Range("F118").Value = "mixing"
$ws.Range("F120").Value = "application instructions"
$ws.Range("F122").Value = "application instructions"
$ws.Range("F124").Value = "application instructions"
$ws.Range("F125").Value = "application instructions"
$ws.Range("F127").Value = "irrigation || application instructions || chemigation"
$ws.Range("F128").Value = "safety procedures"
$ws.Range("F130").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F131").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F132").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F134").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F135").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F136").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F137").Value = "safety procedures || irrigation || application instructions || chemigation"
$ws.Range("F138").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F140").Value = "irrigation || application instructions"
$ws.Range("F141").Value = "irrigation || application instructions"
$ws.Range("F144").Value = "irrigation"
$ws.Range("F146").Value = "irrigation"
$ws.Range("F147").Value = "irrigation"
$ws.Range("F148").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F149").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F152").Value = "safety procedures || irrigation || chemigation"
$ws.Range("F153").Value = "irrigation"
$ws.Range("F556").Value = "154_pesticide_storage"
